$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: push the old row 9 ("Wednesday... / 1/1 / developed by ...") down
# to row 10, preserving its per-cell styles (copy cell-by-cell so the
# engine reuses the existing style indices instead of minting new ones).
# Merge the destination ranges *before* copying values in -- merging
# already-populated cells causes the engine to fork per-cell border
# styles (left/middle/right edge variants); merging blank cells first and
# copying afterwards avoids that.
# ---------------------------------------------------------------------------
$ws.Range("A10:F10").Merge()
$ws.Range("G10:I10").Merge()
$ws.Range("K10:Q10").Merge()

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
foreach ($c in $cols) {
  $ws.Range($c + "9").Copy($ws.Range($c + "10"))
}
$ws.Rows("10:10").RowHeight = 16.5

# remove the old row-9 merges/content now that it lives on row 10
$ws.Range("A9:F9").UnMerge()
$ws.Range("G9:I9").UnMerge()
$ws.Range("K9:Q9").UnMerge()
$ws.Range("A9:Q9").Clear()

# ---------------------------------------------------------------------------
# Step 2: push the old row 8 (just the P8:Q8 total cell) down to row 9.
# ---------------------------------------------------------------------------
# old P8:Q8 merge is no longer valid; row 8 will become a brand-new
# data row shaped like row 7. Unmerge/clear it first...
$ws.Range("P8:Q8").UnMerge()
$ws.Range("A8:Q8").Clear()

# ...then merge the row-9 destination while still blank, and copy in.
$ws.Range("P9:Q9").Merge()
$ws.Range("P8").Copy($ws.Range("P9"))
$ws.Range("Q8").Copy($ws.Range("Q9"))
$ws.Rows("9:9").RowHeight = 25.5

# ---------------------------------------------------------------------------
# Step 3: set the "Text" number format on row 7's text columns BEFORE
# duplicating the row, so the copy in step 4 reuses the same style ids.
# ---------------------------------------------------------------------------
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("Q7").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Step 4: build the new row 8 as a duplicate of row 7's layout/format.
# Merge first (destination still blank), then copy each source cell in.
# ---------------------------------------------------------------------------
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

foreach ($c in $cols) {
  $ws.Range($c + "7").Copy($ws.Range($c + "8"))
}
$ws.Rows("8:8").RowHeight = 24.75

# ---------------------------------------------------------------------------
# Step 5: fill in the actual row 7 / row 8 values.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "DEPO-PEN 1.2 MIU VIAL."
$ws.Range("H7").Value = "3:0"
$ws.Range("L7").Value = "'1"
$ws.Range("N7").Value = "25.00"
$ws.Range("P7").Value = "'25.0000"
$ws.Range("Q7").Value = "1:0"

$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "TORSERETIC 100MG 30 TABS."
$ws.Range("H8").Value = "1:0"
$ws.Range("L8").Value = "'1"
$ws.Range("N8").Value = "261.00"
$ws.Range("P8").Value = "'86.1300"
$ws.Range("Q8").Value = "0:1"

# ---------------------------------------------------------------------------
# Step 6: the new grand-total cell that landed on row 9.
# ---------------------------------------------------------------------------
$ws.Range("P9").Value = 111.13
